$wb = $excel.ActiveWorkbook

# --- 1. Status text change: "Ready for handoff" -> "In Translation" ---
# This shared string shows up on the "Overview" sheet (columns "zh-cn"/"de-de",
# i.e. E2 and F2) as well as on the per-language "zh-cn" and "de-de" sheets
# (column "Status", i.e. C2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Narrow the "zh-cn" / "de-de" status columns ---
# Overview!E:F and the Status column (C) on each language sheet were widened
# from a "Generate Report" layout (~17.22 chars) to a narrower archive layout
# (~13.41 chars). ColumnWidth of 12.5 is the COM value that lands on the
# stored width closest to that target.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
